$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 357 (Royal Glory),
# shifting the existing data (old rows 357-374) down to rows 359-376.
$ws.Rows.Item(357).Insert()
$ws.Rows.Item(357).Insert()

# Fill the two newly inserted rows with the new "Early Majestic" records.
$ws.Range("A357").Value = 11
$ws.Range("B357").Value = "Vega Monumental Concepción"
$ws.Range("C357").Value = "Bíobío"
$ws.Range("D357").Value = 45267
$ws.Range("E357").Value = 8
$ws.Range("F357").Value = "Fruta"
$ws.Range("G357").Value = 100103
$ws.Range("H357").Value = "Frutos de hueso (carozo)"
$ws.Range("I357").Value = 100103004
$ws.Range("J357").Value = "Durazno"
$ws.Range("K357").Value = "Early Majestic"
$ws.Range("L357").Value = "Primera"
$ws.Range("M357").Value = 150
$ws.Range("N357").Value = 15000
$ws.Range("O357").Value = 15000
$ws.Range("P357").Value = 15000
$ws.Range("Q357").Value = "$/caja 15 kilos empedrada"
$ws.Range("R357").Value = "Región de O'Higgins"
$ws.Range("S357").Value = 1000
$ws.Range("T357").Value = 15

$ws.Range("A358").Value = 11
$ws.Range("B358").Value = "Vega Monumental Concepción"
$ws.Range("C358").Value = "Bíobío"
$ws.Range("D358").Value = 45267
$ws.Range("E358").Value = 8
$ws.Range("F358").Value = "Fruta"
$ws.Range("G358").Value = 100103
$ws.Range("H358").Value = "Frutos de hueso (carozo)"
$ws.Range("I358").Value = 100103004
$ws.Range("J358").Value = "Durazno"
$ws.Range("K358").Value = "Early Majestic"
$ws.Range("L358").Value = "Segunda"
$ws.Range("M358").Value = 180
$ws.Range("N358").Value = 13000
$ws.Range("O358").Value = 13000
$ws.Range("P358").Value = 13000
$ws.Range("Q358").Value = "$/caja 15 kilos empedrada"
$ws.Range("R358").Value = "Región de O'Higgins"
$ws.Range("S358").Value = 867
$ws.Range("T358").Value = 15
